$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(2)
Write-Host "Shapes count in footer range:" $ftr.Range.ShapeRange.Count
